$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string text edits (report header text)
#    "Volume 32   Number  27" -> "...28"
# ---------------------------------------------------------------------------
$a8 = $ws.Range("A8")
$a8text = $a8.Value2
$idx = $a8text.LastIndexOf("27")
$a8.Characters($idx + 1, 2).Text = "28"

#    "Report Covering the Week  6/30/2025  Through  7/6/2025"
#      -> "...7/7/2025  Through  7/13/2025"
$c9 = $ws.Range("C9")
$c9text = $c9.Value2
$idx1 = $c9text.IndexOf("6/30/2025")
$c9.Characters($idx1 + 1, "6/30/2025".Length).Text = "7/7/2025"

$c9text2 = $c9.Value2
$idx2 = $c9text2.IndexOf("7/6/2025")
$c9.Characters($idx2 + 1, "7/6/2025".Length).Text = "7/13/2025"

# ---------------------------------------------------------------------------
# 2. Column width tweaks: columns E and H narrow down to match the other
#    numeric columns (same width as column F/G).
# ---------------------------------------------------------------------------
$narrowWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $narrowWidth
$ws.Columns.Item(8).ColumnWidth = $narrowWidth

# ---------------------------------------------------------------------------
# 3. Cells that flip from a shared-string placeholder ("0" / "***.*") to an
#    actual number. Set the value then copy the numeric format from a
#    neighboring numeric cell so the style matches (numFmt "#,##0").
# ---------------------------------------------------------------------------
$numFmtSource = $ws.Range("J14")

$ws.Range("C14").Value = 1
$numFmtSource.Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("C29").Value = 2
$numFmtSource.Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("C30").Value = 2
$numFmtSource.Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Cells that flip from a number to the shared-string placeholder text.
#    Cells receiving the numeric-looking placeholder "0" need the Text
#    number format first so Excel stores it as literal text, not 0.
#    Afterwards copy formatting from a neighboring placeholder-text cell so
#    the final style matches (no numFmt, right aligned).
# ---------------------------------------------------------------------------
$textFmtSource = $ws.Range("D14")

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$textFmtSource.Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").Value = "***.*"
$textFmtSource.Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$textFmtSource.Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$textFmtSource.Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").Value = "***.*"
$textFmtSource.Copy()
$ws.Range("H22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Remaining plain numeric value updates (week/28-day/YTD crime counts and
#    the derived percentage-change figures).
# ---------------------------------------------------------------------------
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -71.428571428571
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -92.307692307692
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 22
$ws.Range("K15").Value = -4.347826086956
$ws.Range("L15").Value = 37.5
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 37.5
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 34
$ws.Range("H16").Value = 54.545454545454
$ws.Range("I16").Value = 140
$ws.Range("J16").Value = 141
$ws.Range("K16").Value = -0.709219858156
$ws.Range("L16").Value = -16.167664670658
$ws.Range("M16").Value = 5.263157894736
$ws.Range("N16").Value = -67.88990825688
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = -27.272727272727
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -17.5
$ws.Range("I17").Value = 239
$ws.Range("J17").Value = 251
$ws.Range("K17").Value = -4.780876494023
$ws.Range("L17").Value = -1.239669421487
$ws.Range("M17").Value = 17.733990147783
$ws.Range("N17").Value = -6.640625
$ws.Range("C18").Value = 7
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -41.379310344827
$ws.Range("I18").Value = 115
$ws.Range("K18").Value = -2.542372881355
$ws.Range("L18").Value = 17.34693877551
$ws.Range("M18").Value = 0.877192982456
$ws.Range("N18").Value = -77.756286266924
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -26.666666666666
$ws.Range("I19").Value = 186
$ws.Range("J19").Value = 224
$ws.Range("K19").Value = -16.964285714285
$ws.Range("L19").Value = -13.888888888888
$ws.Range("M19").Value = 32.857142857142
$ws.Range("N19").Value = -18.777292576419
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 98
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -44.943820224719
$ws.Range("M20").Value = 151.282051282051
$ws.Range("N20").Value = -54.205607476635
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 35.714285714285
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = -9.48905109489
$ws.Range("I21").Value = 802
$ws.Range("J21").Value = 857
$ws.Range("K21").Value = -6.417736289381
$ws.Range("L21").Value = -13.203463203463
$ws.Range("M21").Value = 25.117004680187
$ws.Range("N21").Value = -52.656434474616
$ws.Range("M22").Value = -23.076923076923
$ws.Range("F23").Value = 2
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 35
$ws.Range("I24").Value = 466
$ws.Range("J24").Value = 415
$ws.Range("K24").Value = 12.289156626506
$ws.Range("L24").Value = 9.905660377358
$ws.Range("M24").Value = 57.432432432432
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 21.428571428571
$ws.Range("I25").Value = 117
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = 5.405405405405
$ws.Range("L25").Value = -14.598540145985
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 61
$ws.Range("H26").Value = -4.918032786885
$ws.Range("I26").Value = 332
$ws.Range("J26").Value = 352
$ws.Range("K26").Value = -5.681818181818
$ws.Range("L26").Value = 21.611721611721
$ws.Range("M26").Value = -3.206997084548
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 34
$ws.Range("K27").Value = -10.526315789473
$ws.Range("L27").Value = 13.333333333333
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 59
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = -15.714285714285
$ws.Range("L28").Value = 7.272727272727
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = -30
$ws.Range("L29").Value = -53.333333333333
$ws.Range("M29").Value = -41.666666666666
$ws.Range("N29").Value = -86.538461538461
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = -22.222222222222
$ws.Range("L30").Value = -46.153846153846
$ws.Range("M30").Value = -41.666666666666
$ws.Range("N30").Value = -84.782608695652

